# Error Calculations and Plots
# Applies the edits described by the target diff:
#  - Remove the "RM 232" row and the "SC 92" row (rows shift up, F35 -> F33)
#  - Move several "missing value" markers around (clear some cells that had
#    values, fill in others that were previously blank) to reflect the
#    updated imputation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove rows no longer present in the data (by row label) ---
# Row 26 = "RM 232" in the original layout.
$ws.Rows(26).Delete()
# After the above delete, the row formerly holding "SC 92" is now row 27.
$ws.Rows(27).Delete()

# --- Cell-level value / blank swaps (row numbers below are POST-delete) ---

# Row 5 ("RM 14"): F had 17.66, now blank.
$ws.Range("F5").ClearContents()

# Row 11 ("RM 58"): F was blank, now 17.65.
$ws.Range("F11").Value = 17.65

# Row 19 ("RM 125"): D was blank, now -15.5; F had 17.81, now blank.
$ws.Range("D19").Value = -15.5
$ws.Range("F19").ClearContents()

# Row 21 ("RM 135"): D had -14.3, now blank.
$ws.Range("D21").ClearContents()

# Row 23 ("RM 140"): D was blank, now -13.9.
$ws.Range("D23").Value = -13.9

# Row 25 ("RM 145"): F was blank, now 16.6.
$ws.Range("F25").Value = 16.6

# Row 27 ("SC 101"): D had -14.6, now blank.
$ws.Range("D27").ClearContents()

# Row 29 ("SC 119"): F had 18.06, now blank.
$ws.Range("F29").ClearContents()

# Row 33 ("SC 232"): D was blank, now -14.1.
$ws.Range("D33").Value = -14.1

Write-Output "edits applied"
